# Add 2022-Q3 data
# 1) Insert a new row into the "总计" (summary) sheet for the 2022-Q3 quarter,
#    shifting the existing quarters down by one row.
# 2) Insert a new "2022-Q3" worksheet (cloned from "2022-Q2" so it keeps the
#    same layout/formatting) right after "2022-Q2", populated with the new
#    fund holding details, and rename the quarter sheets don't otherwise change.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet ("总计"): insert the 2022-Q3 row right under the header.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)  # xlPasteFormats -> copy the index-column style

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.65

# The row insert shifted the existing rows down but left their literal
# "index" column (A) values untouched; renumber them 1..7 to match.
for ($i = 3; $i -le 9; $i++) {
    $summary.Range("A$i").Value = $i - 2
}

# ---------------------------------------------------------------------------
# 2. New "2022-Q3" detail worksheet, cloned from "2022-Q2" (index 2) so it
#    inherits the same sheet/page/style setup, placed right before it.
# ---------------------------------------------------------------------------
$quarterTemplate = $wb.Worksheets.Item("2022-Q2")
$quarterTemplate.Copy($quarterTemplate, $null)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# Trim the copied sheet down to just the header + 4 data rows we need.
$newSheet.Rows("6:11").Delete()

$fundRows = @(
    @("516950", "银华中证基建ETF", "11.07", "97.93", "3.60", "0.3985", 9),
    @("159635", "华夏中证基建ETF", "3.40", "99.03", "3.59", "0.1221", 9),
    @("159619", "国泰中证基建ETF", "3.30", "98.76", "3.49", "0.1152", 9),
    @("515870", "嘉实中证先进制造100策略ETF", "0.36", "98.05", "2.65", "0.0095", 7)
)

$r = 2
foreach ($fund in $fundRows) {
    $newSheet.Range("A$r").Value = ($r - 2)

    # Force these as plain text (fund code / name / percentages are text in
    # this workbook, not numbers) via the leading-apostrophe literal marker,
    # then clear the resulting "quote prefix" style so the cell matches the
    # unstyled text cells used elsewhere in the workbook.
    $newSheet.Range("B$r").Value = "'" + $fund[0]
    $newSheet.Range("B$r").ClearFormats()
    $newSheet.Range("C$r").Value = "'" + $fund[1]
    $newSheet.Range("C$r").ClearFormats()
    $newSheet.Range("D$r").Value = "'" + $fund[2]
    $newSheet.Range("D$r").ClearFormats()
    $newSheet.Range("E$r").Value = "'" + $fund[3]
    $newSheet.Range("E$r").ClearFormats()
    $newSheet.Range("F$r").Value = "'" + $fund[4]
    $newSheet.Range("F$r").ClearFormats()
    $newSheet.Range("G$r").Value = "'" + $fund[5]
    $newSheet.Range("G$r").ClearFormats()

    $newSheet.Range("H$r").Value = $fund[6]
    $r++
}

# Restore the original active sheet/selection (the last sheet, "2020-Q4",
# was the active tab before this edit; none of our sheet inserts should
# change that).
$wb.Worksheets.Item("2020-Q4").Select()
